# Updates after WW1916 pilot class in KY.
#
# 1) Refresh the cached "datetimeFigureOut" date field text that appears
#    on every slide layout's Date Placeholder and on the slide master's
#    Date Placeholder (4/10/2019 -> 4/18/2019).
# 2) Bump the version footer textbox on slide 1 (Version 4.0 -> Version 4.1).

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "4/10/2019") {
                $tr.Text = "4/18/2019"
            }
        }
    }
}

# Slide master's Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every custom (slide) layout's Date Placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# Slide 1: "Version 4.0" -> "Version 4.1" footer textbox.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "Version 4.0") {
            $tr.Text = "Version 4.1"
        }
    }
}
